$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row at row 54 (this shifts all the existing
# rows 54..137 down by one, to 55..138, and extends the sheet
# dimension from A1:T137 to A1:T138 automatically).
$ws.Rows(54).EntireRow.Insert()

# Populate the newly inserted row with the new "Elegant Lady" entry.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44579
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103004
$ws.Range("J54").Value = "Durazno"
$ws.Range("K54").Value = "Elegant Lady"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 220
$ws.Range("N54").Value = 13000
$ws.Range("O54").Value = 14000
$ws.Range("P54").Value = 13455
$ws.Range("Q54").Value = "$/bandeja 10 kilos granel"
$ws.Range("R54").Value = "Región de O'Higgins"
$ws.Range("S54").Value = 1346
$ws.Range("T54").Value = 10
